$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 3.5
$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.85
$ws.Range("G3").Value = 1.22
$ws.Range("H3").Value = 5.5
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = 1.73
$ws.Range("K3").Value = 2.4
$ws.Range("L3").Value = 12
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 13
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = 1.83
$ws.Range("R3").Value = 1.98
$ws.Range("S3").Value = 1.36
$ws.Range("T3").Value = 3
$ws.Range("U3").Value = 2.75
$ws.Range("V3").Value = 1.4
$ws.Range("X3").Value = 5
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 6.5
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 9.5
$ws.Range("AD3").Value = 11
$ws.Range("AE3").Value = 34
$ws.Range("AF3").Value = 126
$ws.Range("AI3").Value = 51
$ws.Range("AJ3").Value = 34
$ws.Range("AK3").Value = 201
$ws.Range("AL3").Value = 101
$ws.Range("AM3").Value = 101
$ws.Range("AN3").Value = 2.88
$ws.Range("AO3").Value = 5.5
$ws.Range("AP3").Value = 23
$ws.Range("AQ3").Value = 15
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 3
$ws.Range("AU3").Value = 13
$ws.Range("AV3").Value = 101
$ws.Range("AW3").Value = 12
$ws.Range("AX3").Value = 51
$ws.Range("AY3").Value = 67
$ws.Range("AZ3").Value = 351
$ws.Range("BA3").Value = 351
$ws.Range("G4").Value = 1.75
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 4.5
$ws.Range("J4").Value = 2.4
$ws.Range("L4").Value = 4.75
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 3.5
$ws.Range("Q4").Value = 1.95
$ws.Range("R4").Value = 1.85
$ws.Range("X4").Value = 8
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 17
$ws.Range("AG4").Value = 301
$ws.Range("AH4").Value = 12
$ws.Range("AJ4").Value = 15
$ws.Range("AK4").Value = 51
$ws.Range("AN4").Value = 3.75
$ws.Range("AO4").Value = 9.5
$ws.Range("AQ4").Value = 29
$ws.Range("AU4").Value = 8.5
$ws.Range("AX4").Value = 23
$ws.Range("AY4").Value = 34
$ws.Range("BB4").Value = 251
$ws.Range("H6").Value = 3.5
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("Q6").Value = 1.88
$ws.Range("R6").Value = 1.98
$ws.Range("X6").Value = 12
$ws.Range("Y6").Value = 9.5
$ws.Range("AA6").Value = 19
$ws.Range("AO6").Value = 13
$ws.Range("G7").Value = 3.1
$ws.Range("I7").Value = 2.38
$ws.Range("X7").Value = 15
$ws.Range("AC7").Value = 7
$ws.Range("AN7").Value = 5
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("O8").Value = 1.44
$ws.Range("P8").Value = 2.75
$ws.Range("Q8").Value = 2.3
$ws.Range("R8").Value = 1.62
$ws.Range("G9").Value = 2.4
$ws.Range("H9").Value = 3.3
$ws.Range("J9").Value = 3.1
$ws.Range("K9").Value = 2.1
$ws.Range("L9").Value = 3.5
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 9.5
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.4
$ws.Range("Q9").Value = 2.05
$ws.Range("R9").Value = 1.8
$ws.Range("U9").Value = 1.8
$ws.Range("V9").Value = 1.91
$ws.Range("W9").Value = 8
$ws.Range("X9").Value = 12
$ws.Range("Y9").Value = 9.5
$ws.Range("AB9").Value = 29
$ws.Range("AC9").Value = 9.5
$ws.Range("AD9").Value = 6.5
$ws.Range("AG9").Value = 251
$ws.Range("AH9").Value = 9
$ws.Range("AI9").Value = 15
$ws.Range("AN9").Value = 4.5
$ws.Range("AO9").Value = 13
$ws.Range("AP9").Value = 23
$ws.Range("AQ9").Value = 41
$ws.Range("AS9").Value = 151
$ws.Range("G10").Value = 1.57
$ws.Range("H10").Value = 3.8
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 2.1
$ws.Range("K10").Value = 2.4
$ws.Range("L10").Value = 5
$ws.Range("M10").Value = 1.03
$ws.Range("N10").Value = 15
$ws.Range("O10").Value = 1.18
$ws.Range("P10").Value = 4.5
$ws.Range("Q10").Value = 1.62
$ws.Range("R10").Value = 2.25
$ws.Range("S10").Value = 1.3
$ws.Range("T10").Value = 3.4
$ws.Range("W10").Value = 9
$ws.Range("Z10").Value = 13
$ws.Range("AA10").Value = 12
$ws.Range("AB10").Value = 21
$ws.Range("AC10").Value = 15
$ws.Range("AD10").Value = 8
$ws.Range("AE10").Value = 15
$ws.Range("AH10").Value = 17
$ws.Range("AI10").Value = 29
$ws.Range("AJ10").Value = 17
$ws.Range("AK10").Value = 51
$ws.Range("AL10").Value = 41
$ws.Range("AM10").Value = 41
$ws.Range("AN10").Value = 3.75
$ws.Range("AO10").Value = 8
$ws.Range("AP10").Value = 17
$ws.Range("AQ10").Value = 23
$ws.Range("AR10").Value = 41
$ws.Range("AS10").Value = 101
$ws.Range("AT10").Value = 3.4
$ws.Range("AW10").Value = 7
$ws.Range("AX10").Value = 26
$ws.Range("AZ10").Value = 81
$ws.Range("G12").Value = 3.1
$ws.Range("J12").Value = 3.65
$ws.Range("L12").Value = 2.77
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 7.5
$ws.Range("P12").Value = 2.72
$ws.Range("Q12").Value = 2
$ws.Range("T12").Value = 2.57
$ws.Range("U12").Value = 1.78
$ws.Range("V12").Value = 1.83
$ws.Range("W12").Value = 8.5
$ws.Range("Y12").Value = 11.25
$ws.Range("AA12").Value = 30
$ws.Range("AB12").Value = 40
$ws.Range("AC12").Value = 8.25
$ws.Range("AE12").Value = 14.5
$ws.Range("AF12").Value = 75
$ws.Range("AO12").Value = 17
$ws.Range("AP12").Value = 24
$ws.Range("AR12").Value = 120
$ws.Range("AT12").Value = 2.52
$ws.Range("AW12").Value = 4.1
$ws.Range("AX12").Value = 11.25
$ws.Range("AY12").Value = 18.5
$ws.Range("AZ12").Value = 45
$ws.Range("BA12").Value = 70
$ws.Range("G13").Value = 5.1
$ws.Range("H13").Value = 4.35
$ws.Range("I13").Value = 1.53
$ws.Range("J13").Value = 4.7
$ws.Range("K13").Value = 2.57
$ws.Range("P13").Value = 5.4
$ws.Range("Q13").Value = 1.38
$ws.Range("R13").Value = 2.8
$ws.Range("T13").Value = 3.75
$ws.Range("Y13").Value = 16.5
$ws.Range("AE13").Value = 13
$ws.Range("AH13").Value = 11.75
$ws.Range("AJ13").Value = 8.5
$ws.Range("AN13").Value = 7.4
$ws.Range("AO13").Value = 25
$ws.Range("AP13").Value = 22
$ws.Range("AS13").Value = 175
$ws.Range("AT13").Value = 3.75
